$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 130.11111
$ws.Range("I38").Value = 21.5
$ws.Range("J38").Value = 999
$ws.Range("K38").Value = 64.5
$ws.Range("L38").Value = 2997
$ws.Range("M38").Value = 307.5
$ws.Range("N38").Value = -3741
$ws.Range("H40").Value = 15848.571
$ws.Range("I40").Value = 32124.5
$ws.Range("J40").Value = 12018.941
$ws.Range("K40").Value = 32124.5
$ws.Range("L40").Value = 12018.941
$ws.Range("M40").Value = -31949.5
$ws.Range("N40").Value = -12368.941
$ws.Range("H58").Value = 79
$ws.Range("I58").Value = 79
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 237
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -87
$ws.Range("N58").ClearContents()
$ws.Range("H98").Value = 1335.0278
$ws.Range("I98").Value = 1016.4839
$ws.Range("K98").Value = 1016.4839
$ws.Range("M98").Value = 481.5161000000001
$ws.Range("H122").Value = 1335.0278
$ws.Range("I122").Value = 1016.4839
$ws.Range("K122").Value = 3049.4517
$ws.Range("M122").Value = -599.4516999999996
$ws.Range("H127").Value = 2675.1428
$ws.Range("I127").Value = 1117.0714
$ws.Range("K127").Value = 3351.2142
$ws.Range("M127").Value = 1608.7858
$ws.Range("H129").Value = 1779.7709
$ws.Range("I129").Value = 837.3333
$ws.Range("J129").Value = 1914.4048
$ws.Range("K129").Value = 2511.9999
$ws.Range("L129").Value = 5743.2144
$ws.Range("M129").Value = 2488.0001
$ws.Range("N129").Value = -15743.2144
$ws.Range("H132").Value = 7035.73
$ws.Range("I132").Value = 1404.2115
$ws.Range("J132").Value = 13136.542
$ws.Range("K132").Value = 4212.6345
$ws.Range("L132").Value = 39409.626
$ws.Range("M132").Value = -1682.6345
$ws.Range("N132").Value = -44469.626
$ws.Range("H138").Value = 6526.5557
$ws.Range("J138").Value = 7425.864
$ws.Range("L138").Value = 22277.592
$ws.Range("N138").Value = -32557.592
$ws.Range("H140").Value = 71558.57000000001
$ws.Range("J140").Value = 69192.22
$ws.Range("L140").Value = 69192.22
$ws.Range("N140").Value = -79552.22
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4342.712
$ws.Range("I32").Value = 2456.375
$ws.Range("K32").Value = 2456.375
$ws.Range("M32").Value = -2169.375
$ws.Range("H61").Value = 16142110
$ws.Range("I61").Value = 21749758
$ws.Range("J61").Value = 20119.625
$ws.Range("K61").Value = 21749758
$ws.Range("L61").Value = 20119.625
$ws.Range("M61").Value = -21749546
$ws.Range("N61").Value = -20543.625
$ws.Range("H74").Value = 27779520
$ws.Range("I74").Value = 41668320
$ws.Range("J74").Value = 1916.5
$ws.Range("K74").Value = 41668320
$ws.Range("L74").Value = 1916.5
$ws.Range("M74").Value = -41667446
$ws.Range("N74").Value = -3664.5
$ws.Range("H77").Value = 27779520
$ws.Range("I77").Value = 41668320
$ws.Range("J77").Value = 1916.5
$ws.Range("K77").Value = 208341600
$ws.Range("L77").Value = 9582.5
$ws.Range("M77").Value = -208337232
$ws.Range("N77").Value = -18318.5
$ws.Range("H88").Value = 130750
$ws.Range("I88").Value = 1500
$ws.Range("J88").Value = 173833.33
$ws.Range("K88").Value = 1500
$ws.Range("L88").Value = 173833.33
$ws.Range("M88").Value = -1094
$ws.Range("N88").Value = -174645.33
$ws.Range("H91").Value = 130750
$ws.Range("I91").Value = 1500
$ws.Range("J91").Value = 173833.33
$ws.Range("K91").Value = 1500
$ws.Range("L91").Value = 173833.33
$ws.Range("M91").Value = -96
$ws.Range("N91").Value = -176641.33
$ws.Range("H132").Value = 29786.348
$ws.Range("I132").Value = 34713.06
$ws.Range("K132").Value = 104139.18
$ws.Range("M132").Value = -101609.18
$ws.Range("H133").Value = 71297
$ws.Range("J133").Value = 71297
$ws.Range("L133").Value = 71297
$ws.Range("N133").Value = -76357
$ws.Range("H136").Value = 16142110
$ws.Range("I136").Value = 21749758
$ws.Range("J136").Value = 20119.625
$ws.Range("K136").Value = 65249274
$ws.Range("L136").Value = 60358.875
$ws.Range("M136").Value = -65246724
$ws.Range("N136").Value = -65458.875
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 34247320
$ws.Range("I94").Value = 73385970
$ws.Range("J94").Value = 999.625
$ws.Range("K94").Value = 73385970
$ws.Range("L94").Value = 999.625
$ws.Range("M94").Value = -73385519
$ws.Range("N94").Value = -1901.625
$ws.Range("H107").Value = 27779578
$ws.Range("I107").Value = 62502610
$ws.Range("J107").Value = 1153.5
$ws.Range("K107").Value = 62502610
$ws.Range("L107").Value = 1153.5
$ws.Range("M107").Value = -62500690
$ws.Range("N107").Value = -4993.5
$ws.Range("H134").Value = 2566
$ws.Range("I134").Value = 2125
$ws.Range("J134").Value = 3007
$ws.Range("K134").Value = 6375
$ws.Range("L134").Value = 9021
$ws.Range("M134").Value = -3840
$ws.Range("N134").Value = -14091
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("H140").Value = 181505.42
$ws.Range("J140").Value = 181505.42
$ws.Range("L140").Value = 181505.42
$ws.Range("N140").Value = -191865.42
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 176.46153
$ws.Range("I22").Value = 152.83333
$ws.Range("J22").Value = 196.71428
$ws.Range("K22").Value = 152.83333
$ws.Range("L22").Value = 196.71428
$ws.Range("M22").Value = 197.16667
$ws.Range("N22").Value = -896.71428
$ws.Range("H31").Value = 3268.5186
$ws.Range("I31").Value = 2627.6597
$ws.Range("J31").Value = 7571.4287
$ws.Range("K31").Value = 2627.6597
$ws.Range("L31").Value = 7571.4287
$ws.Range("M31").Value = -2332.6597
$ws.Range("N31").Value = -8161.4287
$ws.Range("H34").Value = 3268.5186
$ws.Range("I34").Value = 2627.6597
$ws.Range("J34").Value = 7571.4287
$ws.Range("K34").Value = 2627.6597
$ws.Range("L34").Value = 7571.4287
$ws.Range("M34").Value = -2425.6597
$ws.Range("N34").Value = -7975.4287
$ws.Range("H58").Value = 1430426
$ws.Range("J58").Value = 1561
$ws.Range("L58").Value = 1561
$ws.Range("N58").Value = -1967
$ws.Range("H99").Value = 12563.318
$ws.Range("J99").Value = 9367.857
$ws.Range("L99").Value = 9367.857
$ws.Range("N99").Value = -12363.857
$ws.Range("H102").Value = 110000
$ws.Range("I102").Value = 100000
$ws.Range("J102").Value = 130000
$ws.Range("K102").Value = 100000
$ws.Range("L102").Value = 130000
$ws.Range("M102").Value = -97566
$ws.Range("N102").Value = -134868
$ws.Range("H122").Value = 3665369.2
$ws.Range("I122").Value = 5918893.5
$ws.Range("J122").Value = 3392.875
$ws.Range("K122").Value = 17756680.5
$ws.Range("L122").Value = 10178.625
$ws.Range("M122").Value = -17754230.5
$ws.Range("N122").Value = -15078.625
$ws.Range("H126").Value = 12563.318
$ws.Range("J126").Value = 9367.857
$ws.Range("L126").Value = 28103.571
$ws.Range("N126").Value = -33043.571
$ws.Range("H132").Value = 2045.5555
$ws.Range("I132").Value = 2058.5715
$ws.Range("K132").Value = 6175.7145
$ws.Range("M132").Value = -3645.7145
$ws.Range("H134").Value = 2626.68
$ws.Range("I134").Value = 2174.8096
$ws.Range("K134").Value = 6524.4288
$ws.Range("M134").Value = -3989.4288
$ws.Range("H136").Value = 1430426
$ws.Range("J136").Value = 1561
$ws.Range("L136").Value = 4683
$ws.Range("N136").Value = -9783
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 786.6667
$ws.Range("J113").Value = 780
$ws.Range("L113").Value = 2340
$ws.Range("N113").Value = -6680
$ws.Range("H131").Value = 2009.9125
$ws.Range("J131").Value = 2048.0908
$ws.Range("L131").Value = 6144.2724
$ws.Range("N131").Value = -16224.2724
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 546.4211
$ws.Range("I97").Value = 655.0833
$ws.Range("K97").Value = 655.0833
$ws.Range("M97").Value = -159.0833
$ws.Range("H113").Value = 866973.75
$ws.Range("I113").Value = 2366215.5
$ws.Range("J113").Value = 10264.286
$ws.Range("K113").Value = 2366215.5
$ws.Range("L113").Value = 10264.286
$ws.Range("M113").Value = -2364045.5
$ws.Range("N113").Value = -14604.286
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7595.6665
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 7595.6665
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 7595.6665
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -7867.6665
$ws.Range("H68").Value = 759664.75
$ws.Range("I68").Value = 1338087.8
$ws.Range("J68").Value = 3265.4614
$ws.Range("K68").Value = 1338087.8
$ws.Range("L68").Value = 3265.4614
$ws.Range("M68").Value = -1337338.8
$ws.Range("N68").Value = -4763.4614
$ws.Range("H71").Value = 759664.75
$ws.Range("I71").Value = 1338087.8
$ws.Range("J71").Value = 3265.4614
$ws.Range("K71").Value = 6690439
$ws.Range("L71").Value = 16327.307
$ws.Range("M71").Value = -6686695
$ws.Range("N71").Value = -23815.307
$ws.Range("H122").Value = 5289.364
$ws.Range("I122").Value = 4934.75
$ws.Range("K122").Value = 14804.25
$ws.Range("M122").Value = -12354.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2779.9246
$ws.Range("I122").Value = 2657.9773
$ws.Range("J122").Value = 3376.111
$ws.Range("K122").Value = 7973.9319
$ws.Range("L122").Value = 10128.333
$ws.Range("M122").Value = -5523.9319
$ws.Range("N122").Value = -15028.333
$ws.Range("H132").Value = 17925534
$ws.Range("I132").Value = 2527349
$ws.Range("K132").Value = 7582047
$ws.Range("M132").Value = -7579517
$ws.Range("H136").Value = 8741.666999999999
$ws.Range("I136").Value = 3819.2307
$ws.Range("J136").Value = 9726.154
$ws.Range("K136").Value = 11457.6921
$ws.Range("L136").Value = 29178.462
$ws.Range("M136").Value = -8907.6921
$ws.Range("N136").Value = -34278.462
